$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[-, 'MEC-3B-Tec. Fundição', -, -]"

$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "[-, -, 'MEC-3B-Tec. Fundição', -]"

$ws.Range("B4").Value = "-"
$ws.Range("D4").Value = "[-, -, 'MEC-3B-Tec. Fundição', -]"
$ws.Range("E4").Value = "['MCT-1A-Tecnologia dos Materiais.', -]"

$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "[-, -, 'MEC-3B-Tec. Fundição', -]"
$ws.Range("E6").Value = "['MCT-1A-Tecnologia dos Materiais.', -]"

$ws.Range("B7").Value = "-"
$ws.Range("F7").Value = "-"

$ws.Range("F8").Value = "-"

$ws.Range("E10").Value = "[-, 'MEC-3A-Tec. Fundição', -, -]"

$ws.Range("E11").Value = "-"

$ws.Range("E12").Value = "-"

$ws.Range("D14").Value = "[-, 'MEC-3A-Tec. Fundição', -, -]"
$ws.Range("E14").Value = "-"

$ws.Range("C15").Value = "[-, -, -, 'MEC-3A-Tec. Fundição']"
$ws.Range("D15").Value = "[-, 'MEC-3A-Tec. Fundição', -, -]"
$ws.Range("E15").Value = "-"

$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("E18").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"

$ws.Range("B19").Value = "[-, -, 'MEC-2NB-Fundição', -]"
$ws.Range("C19").Value = "-"
$ws.Range("E19").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"

$ws.Range("B20").Value = "[-, -, 'MEC-2NB-Fundição', -]"
$ws.Range("C20").Value = "-"
$ws.Range("E20").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"

$ws.Range("B21").Value = "[-, -, 'MEC-2NB-Fundição', -]"
$ws.Range("C21").Value = "[-, 'MEC-2NB-Fundição', -, -]"
$ws.Range("E21").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"
